$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-09-30 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-01 Wednesday", 2) | Out-Null

# Update the division-problem table cells by explicit (row, column)
# position, so we never risk one replacement's new text colliding with
# another replacement's old text (e.g. "23÷7=3, 2" shows up both as a
# target value and, later, as a source value).
$tbl = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "87÷6=14, 3" },
    @{ Row = 1;  Col = 2; New = "40÷6=6, 4" },
    @{ Row = 1;  Col = 3; New = "15÷9=1, 6" },
    @{ Row = 1;  Col = 4; New = "23÷4=5, 3" },
    @{ Row = 1;  Col = 5; New = "12÷9=1, 3" },

    @{ Row = 5;  Col = 1; New = "23÷7=3, 2" },
    @{ Row = 5;  Col = 2; New = "77÷6=12, 5" },
    @{ Row = 5;  Col = 3; New = "37÷2=18, 1" },
    @{ Row = 5;  Col = 4; New = "82÷4=20, 2" },
    @{ Row = 5;  Col = 5; New = "30÷4=7, 2" },

    @{ Row = 9;  Col = 1; New = "64÷6=10, 4" },
    @{ Row = 9;  Col = 2; New = "89÷8=11, 1" },
    @{ Row = 9;  Col = 3; New = "51÷9=5, 6" },
    @{ Row = 9;  Col = 4; New = "14÷8=1, 6" },
    @{ Row = 9;  Col = 5; New = "17÷6=2, 5" },

    @{ Row = 13; Col = 1; New = "84÷9=9, 3" },
    @{ Row = 13; Col = 2; New = "13÷3=4, 1" },
    @{ Row = 13; Col = 3; New = "65÷6=10, 5" },
    @{ Row = 13; Col = 4; New = "52÷2=26, 0" },
    @{ Row = 13; Col = 5; New = "61÷5=12, 1" },

    @{ Row = 17; Col = 1; New = "38÷6=6, 2" },
    @{ Row = 17; Col = 2; New = "10÷3=3, 1" },
    @{ Row = 17; Col = 3; New = "77÷6=12, 5" },
    @{ Row = 17; Col = 4; New = "69÷7=9, 6" },
    @{ Row = 17; Col = 5; New = "15÷6=2, 3" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $cellRange = $cell.Range
    $cellRange.MoveEnd(1, -1) | Out-Null
    $cellRange.Text = $u.New
}

Write-Output "Done"
